$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace all occurrences of the misspelled "califonia cone snail" with the
# corrected "california cone snail" (str_replace_all semantics).
$ws.Cells.Replace("califonia cone snail", "california cone snail")
